# Add new "ACCESS_MODE_CODE" / "EGRESS_MODE_CODE" dictionary rows for the
# Caltrain Pilot operator (progress on tech, oper, paths).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$operator = "Caltrain Pilot"

$rows = @(
    @(562, $operator, "ACCESS_MODE_CODE", 1, "access_mode", "walk"),
    @(563, $operator, "ACCESS_MODE_CODE", 2, "access_mode", "bike"),
    @(564, $operator, "ACCESS_MODE_CODE", 3, "access_mode", "pnr"),
    @(565, $operator, "ACCESS_MODE_CODE", 4, "access_mode", "pnr"),
    @(566, $operator, "ACCESS_MODE_CODE", 5, "access_mode", "pnr"),
    @(567, $operator, "ACCESS_MODE_CODE", 6, "access_mode", "pnr"),
    @(568, $operator, "ACCESS_MODE_CODE", 7, "access_mode", "knr"),
    @(569, $operator, "ACCESS_MODE_CODE", 8, "access_mode", "knr"),
    @(570, $operator, "ACCESS_MODE_CODE", 9, "access_mode", "knr"),
    @(571, $operator, "EGRESS_MODE_CODE", 1, "egress_mode", "walk"),
    @(572, $operator, "EGRESS_MODE_CODE", 2, "egress_mode", "bike"),
    @(573, $operator, "EGRESS_MODE_CODE", 3, "egress_mode", "pnr"),
    @(574, $operator, "EGRESS_MODE_CODE", 4, "egress_mode", "pnr"),
    @(575, $operator, "EGRESS_MODE_CODE", 5, "egress_mode", "pnr"),
    @(576, $operator, "EGRESS_MODE_CODE", 6, "egress_mode", "pnr"),
    @(577, $operator, "EGRESS_MODE_CODE", 7, "egress_mode", "knr"),
    @(578, $operator, "EGRESS_MODE_CODE", 8, "egress_mode", "knr"),
    @(579, $operator, "EGRESS_MODE_CODE", 9, "egress_mode", "knr")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 5).Value2 = $row[5]
}

# Move the active selection to the first empty row below the new data,
# matching where Excel leaves the cursor after entering this block.
$null = $ws.Range("A580").Select()

Write-Output "Added $($rows.Count) dictionary rows (ACCESS_MODE_CODE / EGRESS_MODE_CODE)."
